$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 165
$ws.Range("A165").Value = 7
$ws.Range("B165").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C165").Value = "Ñuble"
$ws.Range("D165").Value = 44939
$ws.Range("D165").NumberFormat = $ws.Range("D164").NumberFormat
$ws.Range("E165").Value = 16
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100103
$ws.Range("H165").Value = "Frutos de hueso (carozo)"
$ws.Range("I165").Value = 100103001
$ws.Range("J165").Value = "Cereza"
$ws.Range("K165").Value = "Sweet Heart"
$ws.Range("L165").Value = "Especial"
$ws.Range("M165").Value = 60
$ws.Range("N165").Value = 6000
$ws.Range("O165").Value = 6000
$ws.Range("P165").Value = 6000
$ws.Range("Q165").Value = "$/bandeja 10 kilos"
$ws.Range("R165").Value = "Quillón"
$ws.Range("S165").Value = 600
$ws.Range("T165").Value = 10

# Row 166
$ws.Range("A166").Value = 7
$ws.Range("B166").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C166").Value = "Ñuble"
$ws.Range("D166").Value = 44939
$ws.Range("D166").NumberFormat = $ws.Range("D164").NumberFormat
$ws.Range("E166").Value = 16
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100103
$ws.Range("H166").Value = "Frutos de hueso (carozo)"
$ws.Range("I166").Value = 100103001
$ws.Range("J166").Value = "Cereza"
$ws.Range("K166").Value = "Sweet Heart"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 60
$ws.Range("N166").Value = 5000
$ws.Range("O166").Value = 5000
$ws.Range("P166").Value = 5000
$ws.Range("Q166").Value = "$/bandeja 10 kilos"
$ws.Range("R166").Value = "Quillón"
$ws.Range("S166").Value = 500
$ws.Range("T166").Value = 10
